# Replace "double" with "var" for the loan/interest/term local variable
# declarations on the two "JavaCalculator01" code slides, per the commit:
# "Changes based on Java in Education meeting - Used var in more places"

$p = $ppt.ActivePresentation

# --- Slide 26: "Classic Java" code block ---
$s26 = $p.Slides.Item(26)
$shp26 = $s26.Shapes.Item(1)
$tr26 = $shp26.TextFrame.TextRange

$tr26.Paragraphs(5, 1).Runs(1).Text = "        var loan = 1000.0;"
$tr26.Paragraphs(6, 1).Runs(1).Text = "`tvar interest = 0.05;"
$tr26.Paragraphs(7, 1).Runs(1).Text = "`tvar term = 5;"

# --- Slide 27: "Nameless Class Java" code block ---
$s27 = $p.Slides.Item(27)
$shp27 = $s27.Shapes.Item(1)
$tr27 = $shp27.TextFrame.TextRange

$tr27.Paragraphs(3, 1).Runs(1).Text = "    var loan = 1000.0;"
$tr27.Paragraphs(4, 1).Runs(1).Text = "    var interest = 0.05;"
$tr27.Paragraphs(5, 1).Runs(1).Text = "    var term = 5;"
